$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.028805385609686
$ws.Range("D2").Value = 1.031653301769922
$ws.Range("E2").Value = 1.038888109898464
$ws.Range("F2").Value = 1.050328011276284
$ws.Range("I2").Value = 1.032163113437563
$ws.Range("J2").Value = 1.033955666876323
$ws.Range("K2").Value = 1.034460797999788
$ws.Range("L2").Value = 1.041674861946162
$ws.Range("M2").Value = 1.053082594528129
$ws.Range("N2").Value = 1.03542400251966

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.029743598430071
$ws.Range("D3").Value = 1.032314048381167
$ws.Range("E3").Value = 1.03974814552809
$ws.Range("F3").Value = 1.051333861952447
$ws.Range("I3").Value = 1.032301477082277
$ws.Range("J3").Value = 1.034534701901018
$ws.Range("K3").Value = 1.034930549589488
$ws.Range("L3").Value = 1.042344871527481
$ws.Range("M3").Value = 1.053900370981181
$ws.Range("N3").Value = 1.036003859840506

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.030351198766547
$ws.Range("D4").Value = 1.032741809770449
$ws.Range("E4").Value = 1.04030548807791
$ws.Range("F4").Value = 1.051985715936248
$ws.Range("I4").Value = 1.032389773096792
$ws.Range("J4").Value = 1.034909283333918
$ws.Range("K4").Value = 1.035234042118263
$ws.Range("L4").Value = 1.042778613177794
$ws.Range("M4").Value = 1.054429915414537
$ws.Range("N4").Value = 1.036378973222006

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.030606755967189
$ws.Range("D5").Value = 1.032921690471237
$ws.Range("E5").Value = 1.040539995117655
$ws.Range("F5").Value = 1.052259993750189
$ws.Range("I5").Value = 1.03242659665088
$ws.Range("J5").Value = 1.035066734129775
$ws.Range("K5").Value = 1.035361517270764
$ws.Range("L5").Value = 1.042961005122156
$ws.Range("M5").Value = 1.054652627670545
$ws.Range("N5").Value = 1.036536647616054

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.030649672271667
$ws.Range("D6").Value = 1.032951896111122
$ws.Range("E6").Value = 1.040579381584214
$ws.Range("F6").Value = 1.052306060143483
$ws.Range("I6").Value = 1.032432762108941
$ws.Range("J6").Value = 1.035093169419834
$ws.Range("K6").Value = 1.035382914247281
$ws.Range("L6").Value = 1.042991632254214
$ws.Range("M6").Value = 1.054690027385166
$ws.Range("N6").Value = 1.036563120447256

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.030354613056173
$ws.Range("D7").Value = 1.032744213151278
$ws.Range("E7").Value = 1.04030862078736
$ws.Range("F7").Value = 1.05198937991473
$ws.Range("I7").Value = 1.032390266299036
$ws.Range("J7").Value = 1.034911387291705
$ws.Range("K7").Value = 1.035235745893212
$ws.Range("L7").Value = 1.042781050124736
$ws.Range("M7").Value = 1.054432890947491
$ws.Range("N7").Value = 1.036381080167654

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.029122352008613
$ws.Range("D8").Value = 1.031876558965625
$ws.Range("E8").Value = 1.039178587819752
$ws.Range("F8").Value = 1.050667735091902
$ws.Range("I8").Value = 1.032210129365901
$ws.Range("J8").Value = 1.034151373125691
$ws.Range("K8").Value = 1.034619649066814
$ws.Range("L8").Value = 1.041901252625066
$ws.Range("M8").Value = 1.053358884767013
$ws.Range("N8").Value = 1.03561998669435

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.026954931163039
$ws.Range("D9").Value = 1.030349345755352
$ws.Range("E9").Value = 1.037193833247526
$ws.Range("F9").Value = 1.048346555567523
$ws.Range("I9").Value = 1.031883274239034
$ws.Range("J9").Value = 1.032811462761969
$ws.Range("K9").Value = 1.033530467551344
$ws.Range("L9").Value = 1.040352521499563
$ws.Range("M9").Value = 1.051469376417369
$ws.Range("N9").Value = 1.034278173504221

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.025512719260002
$ws.Range("D10").Value = 1.029332435645617
$ws.Range("E10").Value = 1.035875118301473
$ws.Range("F10").Value = 1.046804375992544
$ws.Range("I10").Value = 1.031659062305579
$ws.Range("J10").Value = 1.031917797979289
$ws.Range("K10").Value = 1.032802027020297
$ws.Range("L10").Value = 1.039321163122796
$ws.Range("M10").Value = 1.050211810807077
$ws.Range("N10").Value = 1.033383239615045

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.024888886419812
$ws.Range("D11").Value = 1.028892412600554
$ws.Range("E11").Value = 1.03530517327363
$ws.Range("F11").Value = 1.046137862125808
$ws.Range("I11").Value = 1.031560486478385
$ws.Range("J11").Value = 1.03153075052553
$ws.Range("K11").Value = 1.03248606671563
$ws.Range("L11").Value = 1.03887485458
$ws.Range("M11").Value = 1.049667785373072
$ws.Range("N11").Value = 1.032995642509507

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.024657266369088
$ws.Range("D12").Value = 1.0287290158231
$ws.Range("E12").Value = 1.035093631962234
$ws.Range("F12").Value = 1.045890479776449
$ws.Range("I12").Value = 1.031523647598783
$ws.Range("J12").Value = 1.031386972083947
$ws.Range("K12").Value = 1.032368624742203
$ws.Range("L12").Value = 1.038709118285742
$ws.Range("M12").Value = 1.049465787692746
$ws.Range("N12").Value = 1.032851659886046

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.024706945150453
$ws.Range("D13").Value = 1.028764062819692
$ws.Range("E13").Value = 1.035139000955678
$ws.Range("F13").Value = 1.045943535485751
$ws.Range("I13").Value = 1.031531559765837
$ws.Range("J13").Value = 1.031417813565845
$ws.Range("K13").Value = 1.032393820059622
$ws.Range("L13").Value = 1.038744667340377
$ws.Range("M13").Value = 1.049509113351316
$ws.Range("N13").Value = 1.032882545166386

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.024869738604304
$ws.Range("D14").Value = 1.028878905197236
$ws.Range("E14").Value = 1.03528768389769
$ws.Range("F14").Value = 1.046117409528582
$ws.Range("I14").Value = 1.031557445919993
$ws.Range("J14").Value = 1.031518865989633
$ws.Range("K14").Value = 1.032476360556764
$ws.Range("L14").Value = 1.03886115389178
$ws.Range("M14").Value = 1.049651086592352
$ws.Range("N14").Value = 1.032983741096206

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.024970054202303
$ws.Range("D15").Value = 1.02894966969678
$ws.Range("E15").Value = 1.035379313824567
$ws.Range("F15").Value = 1.046224564356224
$ws.Range("I15").Value = 1.031573365642918
$ws.Range("J15").Value = 1.031581126189792
$ws.Range("K15").Value = 1.032527205877285
$ws.Range("L15").Value = 1.038932930772723
$ws.Range("M15").Value = 1.049738571307584
$ws.Range("N15").Value = 1.03304608971299

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.025554134857704
$ws.Range("D16").Value = 1.029361645090474
$ws.Range("E16").Value = 1.035912966273095
$ws.Range("F16").Value = 1.046848636978948
$ws.Range("I16").Value = 1.031665573097654
$ws.Range("J16").Value = 1.031943483339822
$ws.Range("K16").Value = 1.03282298495696
$ws.Range("L16").Value = 1.039350789063781
$ws.Range("M16").Value = 1.0502479268009
$ws.Range("N16").Value = 1.033408961451737

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02592068908032
$ws.Range("D17").Value = 1.029620149243335
$ws.Range("E17").Value = 1.036247999009611
$ws.Range("F17").Value = 1.047240439495974
$ws.Range("I17").Value = 1.03172301369008
$ws.Range("J17").Value = 1.032170758393309
$ws.Range("K17").Value = 1.033008375286796
$ws.Range("L17").Value = 1.039612975273621
$ws.Range("M17").Value = 1.050567569069851
$ws.Range("N17").Value = 1.033636559261868

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.026134557037804
$ws.Range("D18").Value = 1.029770959745524
$ws.Range("E18").Value = 1.036443520822562
$ws.Range("F18").Value = 1.047469092768354
$ws.Range("I18").Value = 1.031756373961848
$ws.Range("J18").Value = 1.032303315850714
$ws.Range("K18").Value = 1.033116458088709
$ws.Range("L18").Value = 1.03976593065674
$ws.Range("M18").Value = 1.050754060009369
$ws.Range("N18").Value = 1.033769304966067

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.026207491156141
$ws.Range("D19").Value = 1.029822387172989
$ws.Range("E19").Value = 1.036510206119549
$ws.Range("F19").Value = 1.047547078194139
$ws.Range("I19").Value = 1.031767724536649
$ws.Range("J19").Value = 1.032348513092468
$ws.Range("K19").Value = 1.033153302639874
$ws.Range("L19").Value = 1.039818088970382
$ws.Range("M19").Value = 1.050817656878687
$ws.Range("N19").Value = 1.033814566393089

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.02588135474981
$ws.Range("D20").Value = 1.029592411142457
$ws.Range("E20").Value = 1.036212042516255
$ws.Range("F20").Value = 1.047198390218115
$ws.Range("I20").Value = 1.031716865733896
$ws.Range("J20").Value = 1.032146374776054
$ws.Range("K20").Value = 1.032988490055258
$ws.Range("L20").Value = 1.039584842415443
$ws.Range("M20").Value = 1.050533269401582
$ws.Range("N20").Value = 1.033612141017081

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.024821797242667
$ws.Range("D21").Value = 1.028845085660133
$ws.Range("E21").Value = 1.035243896004292
$ws.Range("F21").Value = 1.046066202685609
$ws.Range("I21").Value = 1.031549829256366
$ws.Range("J21").Value = 1.031489108887216
$ws.Range("K21").Value = 1.032452056644889
$ws.Range("L21").Value = 1.038826850322244
$ws.Range("M21").Value = 1.049609276860373
$ws.Range("N21").Value = 1.03295394173529

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024156186029383
$ws.Range("D22").Value = 1.028375487267586
$ws.Range("E22").Value = 1.034636119764337
$ws.Range("F22").Value = 1.04535545469847
$ws.Range("I22").Value = 1.03144351436174
$ws.Range("J22").Value = 1.031075792024981
$ws.Range("K22").Value = 1.032114315869521
$ws.Range("L22").Value = 1.038350517423145
$ws.Range("M22").Value = 1.049028775552125
$ws.Range("N22").Value = 1.032540037915703

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024508984348575
$ws.Range("D23").Value = 1.028624403732757
$ws.Range("E23").Value = 1.034958224299585
$ws.Range("F23").Value = 1.04573213065496
$ws.Range("I23").Value = 1.031499996234349
$ws.Range("J23").Value = 1.031294905217198
$ws.Range("K23").Value = 1.032293402275992
$ws.Range("L23").Value = 1.038603006763181
$ws.Range("M23").Value = 1.049336467288901
$ws.Range("N23").Value = 1.032759462273785

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.025899128045785
$ws.Range("D24").Value = 1.029604944705578
$ws.Range("E24").Value = 1.036228289390102
$ws.Range("F24").Value = 1.047217390101504
$ws.Range("I24").Value = 1.031719644175193
$ws.Range("J24").Value = 1.032157392708643
$ws.Range("K24").Value = 1.032997475496952
$ws.Range("L24").Value = 1.039597554360985
$ws.Range("M24").Value = 1.05054876779392
$ws.Range("N24").Value = 1.033623174596397

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.027514783622716
$ws.Range("D25").Value = 1.030743956123078
$ws.Range("E25").Value = 1.037706160499744
$ws.Range("F25").Value = 1.048945712348041
$ws.Range("I25").Value = 1.031968888408209
$ws.Range("J25").Value = 1.033157934920375
$ws.Range("K25").Value = 1.033812460356403
$ws.Range("L25").Value = 1.040752711287918
$ws.Range("M25").Value = 1.051957493512593
$ws.Range("N25").Value = 1.034625137692833
